$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.20"
$ws.Range("E2").Value = "'-6.26%"
$ws.Range("D3").Value = "'35.02"
$ws.Range("E3").Value = "'-3.44%"
$ws.Range("D4").Value = "'5.016"
$ws.Range("E4").Value = "'-1.97%"
$ws.Range("D5").Value = "'0.07901"
$ws.Range("E5").Value = "'-2.48%"
$ws.Range("D6").Value = "'1.935"
$ws.Range("E6").Value = "'-10.06%"
$ws.Range("E7").Value = "'-3.54%"
$ws.Range("D8").Value = "'4.010"
$ws.Range("E8").Value = "'-2.98%"
$ws.Range("D10").Value = "'0.9220"
$ws.Range("E10").Value = "'-0.77%"
$ws.Range("D11").Value = "'0.1165"
$ws.Range("E11").Value = "'16.27%"
$ws.Range("D12").Value = "'0.1834"
$ws.Range("E12").Value = "'-2.56%"
$ws.Range("D13").Value = "'0.09293"
$ws.Range("E13").Value = "'0.73%"
$ws.Range("D14").Value = "'0.03533"
$ws.Range("E14").Value = "'-1.53%"
$ws.Range("D15").Value = "'0.09875"
$ws.Range("E15").Value = "'-0.74%"
$ws.Range("D16").Value = "'0.001384"
$ws.Range("E16").Value = "'-3.78%"
$ws.Range("D17").Value = "'0.005826"
$ws.Range("E17").Value = "'2.41%"
$ws.Range("D18").Value = "'3.496"
$ws.Range("E18").Value = "'0.90%"
$ws.Range("D19").Value = "'0.3443"
$ws.Range("E19").Value = "'2.11%"
$ws.Range("D20").Value = "'0.1308"
$ws.Range("E20").Value = "'-1.65%"
$ws.Range("D21").Value = "'5.037"
$ws.Range("E21").Value = "'-0.26%"
$ws.Range("E22").Value = "'8.87%"
$ws.Range("E23").Value = "'-2.25%"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'-2.37%"
$ws.Range("D25").Value = "'0.004572"
$ws.Range("E25").Value = "'-3.63%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'-3.88%"
$ws.Range("E27").Value = "'-6.87%"
$ws.Range("D39").Value = "'0.01899"
$ws.Range("E39").Value = "'-7.00%"
$ws.Range("D40").Value = "'0.04698"
$ws.Range("D41").Value = "'0.007593"
$ws.Range("E41").Value = "'-2.65%"
$ws.Range("D42").Value = "'0.009540"
$ws.Range("E42").Value = "'22.20%"
$ws.Range("D43").Value = "'0.1323"
$ws.Range("E43").Value = "'-5.56%"
$ws.Range("E44").Value = "'1.35%"
$ws.Range("D45").Value = "'0.01113"
$ws.Range("E45").Value = "'-8.31%"
$ws.Range("D46").Value = "'0.00006002"
$ws.Range("E46").Value = "'-6.32%"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("E49").Value = "'-31.39%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.04%"
